$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (so existing D..H shift to E..I), making room for the new "Status" column
$ws.Range("D1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("A1").Value = "ISIN"
$ws.Range("B1").Value = "Stock Name"
$ws.Range("C1").Value = "Mutual Fund"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Jan_2026"
$ws.Range("F1").Value = "Dec_2025"
$ws.Range("G1").Value = "Oct_2025"
$ws.Range("H1").Value = "MoM"
$ws.Range("I1").Value = "QoQ"

# --- Data rows (A:ISIN, B:Name, C:Fund, D:Status, E:Jan_2026, F:Dec_2025, G:Oct_2025, H:MoM, I:QoQ) ---
$data = @(
    @("INE040A01034","HDFC Bank Limited","quant Aggressive Hybrid Fund","Adding Consistently",8.303181,0.486552,0,7.816629000000001,8.303181),
    @("INE018A01030","Larsen & Toubro Limited","quant Aggressive Hybrid Fund","Adding Consistently",7.967097,7.795835,7.511396,0.1712619999999996,0.4557009999999995),
    @("INE814H01029","Adani Power Limited","quant Aggressive Hybrid Fund","Adding",7.270346,7.235738,7.796691,0.03460800000000042,-0.5263450000000001),
    @("INE917I01010","Bajaj Auto Limited","quant Aggressive Hybrid Fund","Adding Consistently",7.01815,6.437652,5.03351,0.5804980000000004,1.984640000000001),
    @("INE090A01021","ICICI Bank Limited","quant Aggressive Hybrid Fund","Fresh Entry",6.16377,0,0,6.16377,6.16377),
    @("INE364U01010","Adani Green Energy Limited","quant Aggressive Hybrid Fund","Reducing",5.924793,6.64915,3.99776,-0.7243569999999995,1.927033),
    @("INE795G01014","HDFC Life Insurance Co Ltd","quant Aggressive Hybrid Fund","Adding Consistently",5.775625,5.58255,5.318251,0.1930749999999994,0.4573739999999997),
    @("INE423A01024","Adani Enterprises Limited","quant Aggressive Hybrid Fund","Reducing",4.516326,4.717528,3.718617,-0.2012019999999994,0.7977090000000002),
    @("INE775A01035","Samvardhana Motherson International Ltd","quant Aggressive Hybrid Fund","Reducing",3.673839,3.67762,0,-0.003781000000000034,3.673839),
    @("INE406A01037","Aurobindo Pharma Limited","quant Aggressive Hybrid Fund","Adding Consistently",3.444555,3.13663,2.947492,0.307925,0.4970629999999998),
    @("INE758E01017","Jio Financial Services Limited","quant Aggressive Hybrid Fund","Reducing Consistently",3.106554,7.429747,7.543441,-4.323193,-4.436887),
    @("INE127D01025","HDFC Asset Management Company Ltd","quant Aggressive Hybrid Fund","Fresh Entry",2.975697,0,0,2.975697,2.975697),
    @("INE0Z8Z23013","Capital Infra Trust InvIT","quant Aggressive Hybrid Fund","Adding Consistently",1.895889,1.810264,1.814661,0.08562499999999984,0.08122799999999986),
    @("INE721A01047","Shriram Finance Limited","quant Aggressive Hybrid Fund","Fresh Entry",1.450518,0,0,1.450518,1.450518),
    @("INE002A01018","Reliance Industries Limited","quant Aggressive Hybrid Fund","Complete Exit",0,5.255267,7.935728,-5.255267,-7.935728),
    @("INE669C01036","Tech Mahindra Limited","quant Aggressive Hybrid Fund","Complete Exit",0,0.335577,0,-0.335577,0),
    @("INE216A01030","Britannia Industries Limited","quant Aggressive Hybrid Fund","Complete Exit",0,6.408499,6.05352,-6.408499,-6.05352),
    @("INE202B01038","Piramal Finance Ltd","quant Aggressive Hybrid Fund","Complete Exit",0,0,0.114857,0,-0.114857),
    @("INE0J1Y01017","Life Insurance Corporation Of India","quant Aggressive Hybrid Fund","Complete Exit",0,0,6.257899,0,-6.257899),
    @("INE062A01020","State Bank of India","quant Aggressive Hybrid Fund","Complete Exit",0,0,5.010222,0,-5.010222),
    @("INE245A01021","Tata Power Company Limited","quant Aggressive Hybrid Fund","Complete Exit",0,3.840354,5.441064,-3.840354,-5.441064)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $row++
}
